# Commit: "add driver to BasePage constructor some fun with Actions class"
#
# Insert a new worksheet named "move" between "writeReview" and "signIn".
# It gets two cells (A1, A2) both containing the new shared string "driver"
# and becomes the active sheet (selection J20).

$wb = $excel.ActiveWorkbook

$writeReview = $wb.Worksheets.Item("writeReview")

# Add a new worksheet right after "writeReview" (i.e. before "signIn")
$moveSheet = $wb.Worksheets.Add($null, $writeReview)
$moveSheet.Name = "move"

# Populate the new sheet with the "driver" value twice
$moveSheet.Range("A1").Value = "driver"
$moveSheet.Range("A2").Value = "driver"

# Match the saved selection/active cell on the new sheet
$moveSheet.Range("J20").Select()
